$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Mayo de 2020 a las 22:05"

# --- Reorder countries: Sudafrica moved above Colombia/Kuwait (rank change
#     caused Colombia & Kuwait to shift down one row in the sorted table) ---
$ws.Range("A41").Value = "Sudafrica"
$ws.Range("A42").Value = "Colombia"
$ws.Range("A43").Value = "Kuwait"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1522149
$ws.Range("C4").Value = 14376
$ws.Range("E4").Value = 1088702
$ws.Range("G4").Value = 676
$ws.Range("H4").Value = 90789

# --- Row 41: Sudafrica (new stats) ---
$ws.Range("B41").Value = 15515
$ws.Range("C41").Value = 1160
$ws.Range("D41").Value = 7006
$ws.Range("E41").Value = 8245
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = 264

# --- Row 42: Colombia (shifted down from old row 41, values unchanged) ---
$ws.Range("B42").Value = 14939
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 3587
$ws.Range("E42").Value = 10790
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 562

# --- Row 43: Kuwait (shifted down from old row 42, values unchanged) ---
$ws.Range("B43").Value = 14850
$ws.Range("C43").Value = 1048
$ws.Range("D43").Value = 4093
$ws.Range("E43").Value = 10645
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 112

# --- Row 83: Costa de Marfil ---
$ws.Range("B83").Value = 2109
$ws.Range("C83").Value = 48
$ws.Range("D83").Value = 1004
$ws.Range("E83").Value = 1078
$ws.Range("G83").Value = 2
$ws.Range("H83").Value = 27

# --- Row 144: Togo ---
$ws.Range("B144").Value = 301
$ws.Range("C144").Value = 3
$ws.Range("D144").Value = 104
$ws.Range("E144").Value = 186
